$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 586.875
$ws.Range("I2").Value = 282.375
$ws.Range("J2").Value = 739.125
$ws.Range("K2").Value = 282.375
$ws.Range("L2").Value = 739.125
$ws.Range("M2").Value = -169.375
$ws.Range("N2").Value = -965.125
$ws.Range("H9").Value = 1128.7858
$ws.Range("I9").Value = 2116.6667
$ws.Range("K9").Value = 2116.6667
$ws.Range("M9").Value = -1947.6667
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("N18").ClearContents()
$ws.Range("H43").Value = 1160
$ws.Range("I43").Value = 1140
$ws.Range("K43").Value = 1140
$ws.Range("M43").Value = -1071
$ws.Range("H112").Value = 4124.5
$ws.Range("J112").Value = 4999.3335
$ws.Range("L112").Value = 14998.0005
$ws.Range("N112").Value = -17214.0005
$ws.Range("H137").Value = 1528.5
$ws.Range("I137").Value = 1532.3334
$ws.Range("J137").Value = 1522.75
$ws.Range("K137").Value = 4597.0002
$ws.Range("L137").Value = 4568.25
$ws.Range("M137").Value = -2047.0002
$ws.Range("N137").Value = -9668.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1674.75
$ws.Range("I45").Value = 1399.6666
$ws.Range("J45").Value = 2500
$ws.Range("K45").Value = 1399.6666
$ws.Range("L45").Value = 2500
$ws.Range("M45").Value = -1022.6666
$ws.Range("N45").Value = -3254
$ws.Range("H74").Value = 2630.5
$ws.Range("I74").Value = 2630.5
$ws.Range("K74").Value = 2630.5
$ws.Range("M74").Value = -1756.5
$ws.Range("H77").Value = 2630.5
$ws.Range("I77").Value = 2630.5
$ws.Range("K77").Value = 13152.5
$ws.Range("M77").Value = -8784.5
$ws.Range("H132").Value = 3367.5557
$ws.Range("I132").Value = 2663.5
$ws.Range("K132").Value = 7990.5
$ws.Range("M132").Value = -5460.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 4318.6665
$ws.Range("I82").Value = 4318.6665
$ws.Range("K82").Value = 4318.6665
$ws.Range("M82").Value = -3935.6665
$ws.Range("H85").Value = 4318.6665
$ws.Range("I85").Value = 4318.6665
$ws.Range("K85").Value = 4318.6665
$ws.Range("M85").Value = -2992.6665
$ws.Range("H94").Value = 795.2857
$ws.Range("I94").Value = 912.25
$ws.Range("K94").Value = 912.25
$ws.Range("M94").Value = -461.25
$ws.Range("H105").Value = 1899
$ws.Range("I105").Value = 1899
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1899
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -152
$ws.Range("N105").ClearContents()
$ws.Range("H107").Value = 474.0909
$ws.Range("I107").Value = 203.07143
$ws.Range("K107").Value = 203.07143
$ws.Range("M107").Value = 1716.92857

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3809.2
$ws.Range("I16").Value = 2265.889
$ws.Range("K16").Value = 2265.889
$ws.Range("M16").Value = -1978.889
$ws.Range("H58").Value = 2264.6
$ws.Range("J58").Value = 2661.5
$ws.Range("L58").Value = 2661.5
$ws.Range("N58").Value = -3067.5
$ws.Range("H113").Value = 3809.2
$ws.Range("I113").Value = 2265.889
$ws.Range("K113").Value = 2265.889
$ws.Range("M113").Value = -95.88900000000012
$ws.Range("H136").Value = 2264.6
$ws.Range("J136").Value = 2661.5
$ws.Range("L136").Value = 7984.5
$ws.Range("N136").Value = -13084.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 7807.5386
$ws.Range("I6").Value = 8457.666999999999
$ws.Range("J6").Value = 6
$ws.Range("K6").Value = 25373.001
$ws.Range("L6").Value = 18
$ws.Range("M6").Value = -25260.001
$ws.Range("N6").Value = -244
$ws.Range("H75").Value = 7000
$ws.Range("H78").Value = 7000

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4335.75
$ws.Range("I80").Value = 2040.25
$ws.Range("J80").Value = 6631.25
$ws.Range("K80").Value = 2040.25
$ws.Range("L80").Value = 6631.25
$ws.Range("M80").Value = -1042.25
$ws.Range("N80").Value = -8627.25
$ws.Range("H83").Value = 4335.75
$ws.Range("I83").Value = 2040.25
$ws.Range("J83").Value = 6631.25
$ws.Range("K83").Value = 10201.25
$ws.Range("L83").Value = 33156.25
$ws.Range("M83").Value = -5209.25
$ws.Range("N83").Value = -43140.25
$ws.Range("H102").Value = 2149.6428
$ws.Range("I102").Value = 2144.077
$ws.Range("K102").Value = 2144.077
$ws.Range("M102").Value = -522.0770000000002
$ws.Range("H113").Value = 1325
$ws.Range("J113").Value = 1325
$ws.Range("L113").Value = 1325
$ws.Range("N113").Value = -5665
$ws.Range("H132").Value = 19989.5
$ws.Range("I132").Value = 19989.5
$ws.Range("K132").Value = 59968.5
$ws.Range("M132").Value = -57438.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 13743.5
$ws.Range("J18").Value = 16691.334
$ws.Range("L18").Value = 16691.334
$ws.Range("N18").Value = -17035.334
$ws.Range("H25").Value = 299
$ws.Range("I25").Value = 299
$ws.Range("K25").Value = 299
$ws.Range("M25").Value = -69
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("N71").ClearContents()
$ws.Range("H93").Value = 1099.3334
$ws.Range("I93").Value = 899
$ws.Range("J93").Value = 1500
$ws.Range("K93").Value = 899
$ws.Range("L93").Value = 1500
$ws.Range("M93").Value = 349
$ws.Range("N93").Value = -3996
$ws.Range("H136").Value = 4108
$ws.Range("I136").Value = 4128.6
$ws.Range("K136").Value = 12385.8
$ws.Range("M136").Value = -9835.800000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 3407
$ws.Range("I6").Value = 3498.6667
$ws.Range("J6").Value = 3367.7144
$ws.Range("K6").Value = 3498.6667
$ws.Range("L6").Value = 3367.7144
$ws.Range("M6").Value = -3383.6667
$ws.Range("N6").Value = -3597.7144
$ws.Range("H9").Value = 1664.6666
$ws.Range("J9").Value = 1664.6666
$ws.Range("L9").Value = 1664.6666
$ws.Range("N9").Value = -1944.6666
$ws.Range("H41").Value = 19699.666
$ws.Range("I41").Value = 19675.5
$ws.Range("J41").Value = 19711.75
$ws.Range("K41").Value = 19675.5
$ws.Range("L41").Value = 19711.75
$ws.Range("M41").Value = -19285.5
$ws.Range("N41").Value = -20491.75
$ws.Range("H81").Value = 3336666.2
$ws.Range("I81").Value = 4998
$ws.Range("K81").Value = 9996
$ws.Range("M81").Value = -8935
$ws.Range("H84").Value = 3336666.2
$ws.Range("I84").Value = 4998
$ws.Range("K84").Value = 49980
$ws.Range("M84").Value = -44676
$ws.Range("H126").Value = 3471.25
$ws.Range("I126").Value = 3249
$ws.Range("K126").Value = 9747
$ws.Range("M126").Value = -7277

